# Apply weekly report refresh: updated generation timestamp, populated
# billed amount / pricing figures now that the report reflects completed
# billing data, and cleared the (now-unused) Scope ID value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Report generated timestamp (D5)
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:00 AM"

# Total Billed Amount (C8)
$ws.Range("C8").Value = 7178.25

# Scope ID # value cleared (G10)
$ws.Range("G10").Value = ""

# Tuesday (07/22/2025) section - single line item + total
$ws.Range("H16").Value = 478.55
$ws.Range("H17").Value = 478.55

# Sunday (07/27/2025) section - fourteen line items, each priced the same
$sundayRows = 22..35
foreach ($r in $sundayRows) {
    $ws.Cells.Item($r, 8).Value = 478.55
}

# Sunday section TOTAL
$ws.Range("H36").Value = 6699.700000000002
